# Apply updated crypto price/volume figures (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: column D holds plain-text numbers (e.g. "491.40"). Several new
# values look like valid numbers to Excel's auto-detection, which would
# silently convert them to numeric cells (losing the exact text, e.g.
# "0.840" -> 0.84). A leading apostrophe forces Excel to keep them as text,
# exactly like a user typing an apostrophe before a numeric-looking entry.

$ws.Range("D2").Value = "53.928.19"
$ws.Range("D3").Value = "2.262.60"
$ws.Range("E3").Value = "  -4.57%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'491.40"
$ws.Range("E5").Value = "  -2.97%  "
$ws.Range("D6").Value = "'126.81"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.526"
$ws.Range("E8").Value = "  -3.24%  "
$ws.Range("D9").Value = "2.260.73"
$ws.Range("E9").Value = "  -5.03%  "
$ws.Range("D10").Value = "'0.0933"
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "'0.322"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("E13").Value = "  -4.92%  "
$ws.Range("D14").Value = "2.653.30"
$ws.Range("E14").Value = "  -4.98%  "
$ws.Range("D15").Value = "'21.48"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "53.874.93"
$ws.Range("E16").Value = "  -4.25%  "
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("D18").Value = "2.253.18"
$ws.Range("E18").Value = "  -5.22%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'9.76"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.03"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "'297.81"
$ws.Range("E21").Value = "  -3.62%  "
$ws.Range("D22").Value = "'6.27"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "'63.84"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").Value = "2.327.52"
$ws.Range("E28").Value = "  -6.28%  "
$ws.Range("D29").Value = "'7.11"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("D30").Value = "'163.11"
$ws.Range("E30").Value = "  -5.76%  "
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("D32").Value = "0.0₃0679"
$ws.Range("E32").Value = "  -4.41%  "
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").Value = "'0.995"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "'17.42"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").Value = "'0.840"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("D40").Value = "'3.61"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").Value = "'35.36"
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("D42").Value = "'0.373"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "'3.33"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "'4.82"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").Value = "'0.0891"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "'242.49"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("D50").Value = "'0.0478"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "'0.0203"
$ws.Range("E51").Value = "  -1.95%  "
